# "Atualizacao de bases das ligas, do dia: 28-05-2024 as 19:13"
#
# The league database refresh corrected several mis-matched fixture
# records: for a handful of row-pairs (and one row-triple) the data
# payload (id + HomeTeam..PL_AhUnder, i.e. columns B and E:AD) had been
# attached to the wrong fixture. Column A (row #), C (Div) and D (Date)
# are correct already and stay put; only the B/E:AD payload needs to move
# between rows.
#
# This script reads each row's B and E:AD payload into memory first, then
# writes the corrected payloads back, so every row in a group is updated
# from a consistent "before" snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowId($r) {
    return $ws.Range("B$r").Value()
}

function Get-RowRest($r) {
    return $ws.Range("E$r" + ":AD$r").Value()
}

function Set-RowId($r, $v) {
    $ws.Range("B$r").Value = $v
}

function Set-RowRest($r, $v) {
    $ws.Range("E$r" + ":AD$r").Value = $v
}

function Swap-RowPayload($r1, $r2) {
    $id1 = Get-RowId $r1
    $rest1 = Get-RowRest $r1
    $id2 = Get-RowId $r2
    $rest2 = Get-RowRest $r2

    Set-RowId $r1 $id2
    Set-RowRest $r1 $rest2

    Set-RowId $r2 $id1
    Set-RowRest $r2 $rest1
}

# Straight two-row swaps.
Swap-RowPayload 29 30
Swap-RowPayload 111 112
Swap-RowPayload 122 123
Swap-RowPayload 198 199

# Rows 189 / 192 / 194 rotate three ways:
#   row189 <- old row194, row192 <- old row189, row194 <- old row192
$id189 = Get-RowId 189
$rest189 = Get-RowRest 189
$id192 = Get-RowId 192
$rest192 = Get-RowRest 192
$id194 = Get-RowId 194
$rest194 = Get-RowRest 194

Set-RowId 189 $id194
Set-RowRest 189 $rest194

Set-RowId 192 $id189
Set-RowRest 192 $rest189

Set-RowId 194 $id192
Set-RowRest 194 $rest192
